$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: convert existing text values to real numbers where appropriate ---
$ws.Range("B2").Value = 1000
$ws.Range("C2").Value = 8000
$ws.Range("D2").Value = 315
$ws.Range("E2").Value = 26
# F2 ("All") stays text - untouched
$ws.Range("G2").Value = 25
$ws.Range("H2").Value = 27
$ws.Range("I2").Value = 26
$ws.Range("J2").Value = 600
$ws.Range("K2").Value = 400
$ws.Range("L2").Value = 1000
$ws.Range("M2").Value = 4
$ws.Range("N2").Value = 96
# O2, P2, Q2 (list-like text) stay text - untouched

# --- Row 3: brand-new scholarship row ---
$ws.Range("A3").Value = "New Scholarship"
# B3 / C3 look numeric but must be stored as TEXT - leading apostrophe forces text entry,
# then reset Style back to Normal so no stray number-format/quote-prefix style lingers on the cell
$ws.Range("B3").Value = "'10"
$ws.Range("B3").Style = "Normal"
$ws.Range("C3").Value = "'1000"
$ws.Range("C3").Style = "Normal"
$ws.Range("D3").Value = 210
$ws.Range("E3").Value = 23
$ws.Range("F3").Value = "Computer Science and Engineering"
$ws.Range("G3").Value = 23
$ws.Range("H3").Value = 19
$ws.Range("I3").Value = 21
$ws.Range("J3").Value = 530
$ws.Range("K3").Value = 370
$ws.Range("L3").Value = 900
$ws.Range("M3").Value = 3.1
$ws.Range("N3").Value = 75
$ws.Range("O3").Value = "['ACT Composite', 'SAT Combined']"
$ws.Range("P3").Value = "['ACT Math', 'SAT Math']"
$ws.Range("Q3").Value = "['GPA', 'HS Percentile']"
